$d = $word.ActiveDocument

# 1. Title: "Data Aware Selection Tools" -> "Density-Aware Selection Tools"
$r = $d.Content
$found = $r.Find.Execute("Data Aware Selection Tools", $true, $false, $false, $false, $false, $true, 1, $false, "Density-Aware Selection Tools", 2)

# 2. Abstract: "We provide data-aware interactions" -> "We provide density-aware interactions"
$r = $d.Content
$found = $r.Find.Execute("We provide data-aware interactions", $true, $false, $false, $false, $false, $true, 1, $false, "We provide density-aware interactions", 2)

# 3. Abstract: "our proposed data-aware slider" -> "our proposed density-aware slider"
$r = $d.Content
$found = $r.Find.Execute("our proposed data-aware slider", $true, $false, $false, $false, $false, $true, 1, $false, "our proposed density-aware slider", 2)

# 4. Abstract: "a data-aware range slider" -> "a density-aware range slider"
$r = $d.Content
$found = $r.Find.Execute("a data-aware range slider", $true, $false, $false, $false, $false, $true, 1, $false, "a density-aware range slider", 2)

# 5. "a given pixels maps few or many items." -> "a given pixels maps to few or many items."
$r = $d.Content
$found = $r.Find.Execute("a given pixels maps few or many items.", $true, $false, $false, $false, $false, $true, 1, $false, "a given pixels maps to few or many items.", 2)

# 6. "propose density aware selection tools" -> "propose density-aware selection tools"
$r = $d.Content
$found = $r.Find.Execute("propose density aware selection tools", $true, $false, $false, $false, $false, $true, 1, $false, "propose density-aware selection tools", 2)

# 7. Move the _GoBack bookmark from the end of the "We propose...querying. " paragraph
#    to just after "To solve the issue of data " in the following paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$r = $d.Content
$found = $r.Find.Execute("To solve the issue of data ", $true, $false, $false, $false, $false, $true, 1, $false, "", $null)
$target = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $target)

Write-Output "done"
